$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'325.51"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-3.63%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'44.51"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'0.75%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.581"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-3.30%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08063"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-3.30%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'8.675"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-1.93%"
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'-4.79%"
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'-3.65%"
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'-6.83%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.9427"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-0.30%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1167"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-6.32%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.1860"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-4.88%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.09946"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'1.18%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.04269"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-5.64%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.1067"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'0.001292"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-0.77%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.04190"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-5.45%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.005975"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-1.77%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'3.587"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'2.55%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.3502"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-0.27%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'8.367"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-4.64%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.1370"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-0.17%"
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'-6.21%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001243"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-1.55%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004471"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'2.34%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001182"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-6.31%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0003995"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'-0.01%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.02621"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-6.46%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.05441"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-6.50%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007690"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-3.26%"
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'-2.54%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.007131"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-20.50%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002019"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-4.89%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.008564"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-4.36%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00007132"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-2.18%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-0.01%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.003659"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'14.51%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.002273"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-0.01%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002103"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.01%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002003"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'-0.01%"
$ws.Range("E51").Style = "Normal"
